$wb = $excel.ActiveWorkbook

# --- Overview sheet: rows 2-5 get reordered so that 068c0ebe moves to the
#     top (its status flips from "Ready for handoff" to "In Translation"),
#     06aa7541 and ee8ee80d shift down one row each, 8a9aaca5 stays last. ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("A2").Value = "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md"
$ovw.Range("B2").Value = "In Translation"
$ovw.Range("C2").Value = "In Translation"

$ovw.Range("A3").Value = "06aa7541-cd06-465c-8316-7632d9c3aa5a.md"
$ovw.Range("B3").Value = "In Translation"
$ovw.Range("C3").Value = "In Translation"

$ovw.Range("A4").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$ovw.Range("B4").Value = "In Translation"
$ovw.Range("C4").Value = "In Translation"

$ovw.Range("A5").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$ovw.Range("B5").Value = "Ready for handoff"
$ovw.Range("C5").Value = "Ready for handoff"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md"
$zh.Range("B2").Value = "In Translation"
$zh.Range("C2").Value = "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.zh-cn.xlf"
$zh.Range("D2").Value = "2016-02-17 09:18:18"

$zh.Range("A3").Value = "06aa7541-cd06-465c-8316-7632d9c3aa5a.md"
$zh.Range("B3").Value = "In Translation"
$zh.Range("C3").Value = "06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.zh-cn.xlf"
$zh.Range("D3").Value = "2016-02-17 09:16:42"

$zh.Range("A4").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$zh.Range("B4").Value = "In Translation"
$zh.Range("C4").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.zh-cn.xlf"
$zh.Range("D4").Value = "2016-02-17 09:16:42"

$zh.Range("A5").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("C5").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.zh-cn.xlf"
$zh.Range("D5").Value = "2016-02-17 09:19:02"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md"
$de.Range("B2").Value = "In Translation"
$de.Range("C2").Value = "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.de-de.xlf"
$de.Range("D2").Value = "2016-02-17 09:18:29"

$de.Range("A3").Value = "06aa7541-cd06-465c-8316-7632d9c3aa5a.md"
$de.Range("B3").Value = "In Translation"
$de.Range("C3").Value = "06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.de-de.xlf"
$de.Range("D3").Value = "2016-02-17 09:16:56"

$de.Range("A4").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$de.Range("B4").Value = "In Translation"
$de.Range("C4").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.de-de.xlf"
$de.Range("D4").Value = "2016-02-17 09:16:56"

$de.Range("A5").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$de.Range("B5").Value = "Ready for handoff"
$de.Range("C5").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.de-de.xlf"
$de.Range("D5").Value = "2016-02-17 09:19:13"
